# Update countries & provincias Spain
# This script applies the latest COVID-19 data refresh to the "Pais" sheet:
#  - A handful of rows swap the displayed country name with their neighbour
#    (this mirrors the reordering that happened in the workbook's shared
#    string table), and
#  - updated case statistics (B..H) for the affected rows, and
#  - a refreshed "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refresh "last updated" timestamp -----------------------------
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 12:22"

# --- Row 4: Estados Unidos --------------------------------------------------
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 7724207
$ws.Range("C4").Value = 1461
$ws.Range("D4").Value = 4936501
$ws.Range("E4").Value = 2571857
$ws.Range("G4").Value = 27
$ws.Range("H4").Value = 215849

# --- Row 19: Banglades ------------------------------------------------------
$ws.Range("A19").Value = "Banglades"
$ws.Range("B19").Value = 373151
$ws.Range("C19").Value = 1520
$ws.Range("D19").Value = 286631
$ws.Range("E19").Value = 81080
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 5440

# --- Rows 31-32: Ecuador / Rumania swap places ------------------------------
$ws.Range("A31").Value = "Rumania"
$ws.Range("B31").Value = 142570
$ws.Range("C31").Value = 2958
$ws.Range("D31").Value = 111564
$ws.Range("E31").Value = 25803
$ws.Range("G31").Value = 82
$ws.Range("H31").Value = 5203

$ws.Range("A32").Value = "Ecuador"
$ws.Range("B32").Value = 142056
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 120511
$ws.Range("E32").Value = 9843
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 11702

# --- Row 43: Oman ------------------------------------------------------------
$ws.Range("A43").Value = "Oman"
$ws.Range("B43").Value = 103465
$ws.Range("C43").Value = 817
$ws.Range("D43").Value = 91329
$ws.Range("E43").Value = 11136
$ws.Range("G43").Value = 10
$ws.Range("H43").Value = 1000

# --- Rows 96-97: Albania / Eslovaquia swap places ---------------------------
$ws.Range("A96").Value = "Eslovaquia"
$ws.Range("B96").Value = 14689
$ws.Range("C96").Value = 877
$ws.Range("D96").Value = 5200
$ws.Range("E96").Value = 9434
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 55

$ws.Range("A97").Value = "Albania"
$ws.Range("B97").Value = 14568
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 8965
$ws.Range("E97").Value = 5200
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 403

# --- Rows 98-99: Sudan / Malasia swap places --------------------------------
$ws.Range("A98").Value = "Malasia"
$ws.Range("B98").Value = 13993
$ws.Range("C98").Value = 489
$ws.Range("D98").Value = 10501
$ws.Range("E98").Value = 3351
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 141

$ws.Range("A99").Value = "Sudan"
$ws.Range("B99").Value = 13653
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 6764
$ws.Range("E99").Value = 6053
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 836

# --- Row 102: Finlandia -------------------------------------------------------
$ws.Range("A102").Value = "Finlandia"
$ws.Range("B102").Value = 11049
$ws.Range("C102").Value = 120
$ws.Range("D102").Value = 8100
$ws.Range("E102").Value = 2603
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 346

# --- Rows 116-117: Jamaica / Eslovenia swap places --------------------------
$ws.Range("A116").Value = "Eslovenia"
$ws.Range("B116").Value = 7120
$ws.Range("C116").Value = 356
$ws.Range("D116").Value = 4535
$ws.Range("E116").Value = 2426
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 159

$ws.Range("A117").Value = "Jamaica"
$ws.Range("B117").Value = 7109
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 2674
$ws.Range("E117").Value = 4312
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 123

# --- Rows 124-125: Republica de Yibuti / Lituania swap places ---------------
$ws.Range("A124").Value = "Lituania"
$ws.Range("B124").Value = 5483
$ws.Range("C124").Value = 117
$ws.Range("D124").Value = 2600
$ws.Range("E124").Value = 2782
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 101

$ws.Range("A125").Value = "Republica de Yibuti"
$ws.Range("B125").Value = 5423
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 5353
$ws.Range("E125").Value = 9
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 61

# --- Row 181: Gibraltar -------------------------------------------------------
$ws.Range("A181").Value = "Gibraltar"
$ws.Range("B181").Value = 445
$ws.Range("C181").Value = 8
$ws.Range("D181").Value = 376
$ws.Range("E181").Value = 69
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 0

# --- Rows 207-208: Nueva Caledonia / Santa Lucia swap places ----------------
# (case counts are identical for both countries, so only the names swap)
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Nueva Caledonia"
